$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new BOM line (row 22) for the M.2 connector J6
$ws.Range("A22").Value = 1
$ws.Range("C22").Value = "A115899CT-ND"
$ws.Range("B22").Value = "J6"

# Match the wrapText style used by the other Digikey_PN cells in column C
$ws.Range("C22").WrapText = $true

# Update the view to match the author's saved selection/scroll state
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E17").Select()
